$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit performs a 3-way cyclic rotation of the observation records
# currently sitting in rows 5, 6 and 7 of the sheet:
#   - the data that was in row 6 moves up into row 5
#   - the data that was in row 7 moves up into row 6
#   - the data that was in row 5 moves down into row 7
# Only the cells actually touched by the source edit are written here.

# ---- Row 5 (becomes what used to be row 6's record) ----
$ws.Range("A5").Value = 111891400
$ws.Range("B5").Value = 96348
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = "Knärot"
$ws.Range("G5").Value = "Goodyera repens"
$ws.Range("H5").Value = "(L.) R. Br."
$ws.Range("I5").Value = "'10"
$ws.Range("J5").Value = "plantor/tuvor"
$ws.Range("K5").Value = "fullt utvecklade blad"
$ws.Range("Q5").Value = 575106.7474690104
$ws.Range("R5").Value = 6702828.414323498
$ws.Range("Z5").Value = "14:38"
$ws.Range("AB5").Value = "14:38"

# ---- Row 6 (becomes what used to be row 7's record) ----
$ws.Range("A6").Value = 111891126
$ws.Range("I6").Value = "'5"
$ws.Range("Q6").Value = 575125.6693508058
$ws.Range("R6").Value = 6702845.585682714
$ws.Range("Z6").Value = "14:33"
$ws.Range("AB6").Value = "14:33"

# ---- Row 7 (becomes what used to be row 5's record) ----
$ws.Range("A7").Value = 111892438
$ws.Range("B7").Value = 55395
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 208257
$ws.Range("F7").Value = "Kopparödla"
$ws.Range("G7").Value = "Anguis fragilis"
$ws.Range("H7").Value = "Linnaeus, 1758"
$ws.Range("I7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("Q7").Value = 575076.7844513921
$ws.Range("R7").Value = 6702914.648038276
$ws.Range("Z7").Value = "15:00"
$ws.Range("AB7").Value = "15:00"
